$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.025004462179452
$ws.Cells.Item(2, 4).Value = 1.033247798436321
$ws.Cells.Item(2, 5).Value = 1.028631126780108
$ws.Cells.Item(2, 6).Value = 1.040182332634911
$ws.Cells.Item(2, 9).Value = 1.031639531406335
$ws.Cells.Item(2, 10).Value = 1.030175907542011
$ws.Cells.Item(2, 11).Value = 1.036050695691546
$ws.Cells.Item(2, 12).Value = 1.03144738233814
$ws.Cells.Item(2, 13).Value = 1.042965406825756
$ws.Cells.Item(2, 14).Value = 1.014093175471355

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025941250642918
$ws.Cells.Item(3, 4).Value = 1.0339643505983
$ws.Cells.Item(3, 5).Value = 1.02951522696082
$ws.Cells.Item(3, 6).Value = 1.041091116857807
$ws.Cells.Item(3, 9).Value = 1.031802521816234
$ws.Cells.Item(3, 10).Value = 1.030751910269369
$ws.Cells.Item(3, 11).Value = 1.036576435339062
$ws.Cells.Item(3, 12).Value = 1.032139253270252
$ws.Cells.Item(3, 13).Value = 1.043684303002206
$ws.Cells.Item(3, 14).Value = 1.014284739735747

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.026547693437132
$ws.Cells.Item(4, 4).Value = 1.034427837185092
$ws.Cells.Item(4, 5).Value = 1.030087932679007
$ws.Cells.Item(4, 6).Value = 1.041679525126929
$ws.Cells.Item(4, 9).Value = 1.031906136428874
$ws.Cells.Item(4, 10).Value = 1.03112430196314
$ws.Cells.Item(4, 11).Value = 1.036915784234094
$ws.Cells.Item(4, 12).Value = 1.032586940908179
$ws.Cells.Item(4, 13).Value = 1.044149211231412
$ws.Cells.Item(4, 14).Value = 1.014408549014147

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.026802707625267
$ws.Cells.Item(5, 4).Value = 1.03462264466745
$ws.Cells.Item(5, 5).Value = 1.030328848604258
$ws.Cells.Item(5, 6).Value = 1.041926977717732
$ws.Cells.Item(5, 9).Value = 1.031949252117797
$ws.Cells.Item(5, 10).Value = 1.031280777731725
$ws.Cells.Item(5, 11).Value = 1.037058244342417
$ws.Cells.Item(5, 12).Value = 1.032775147917316
$ws.Cells.Item(5, 13).Value = 1.044344593796908
$ws.Cells.Item(5, 14).Value = 1.014460563220377

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.026845529451834
$ws.Cells.Item(6, 4).Value = 1.034655351179621
$ws.Cells.Item(6, 5).Value = 1.030369308232651
$ws.Cells.Item(6, 6).Value = 1.041968531084067
$ws.Cells.Item(6, 9).Value = 1.03195646537986
$ws.Cells.Item(6, 10).Value = 1.031307046117775
$ws.Cells.Item(6, 11).Value = 1.037082152133773
$ws.Cells.Item(6, 12).Value = 1.032806748619465
$ws.Cells.Item(6, 13).Value = 1.044377395549799
$ws.Cells.Item(6, 14).Value = 1.014469294559141

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.026551100692438
$ws.Cells.Item(7, 4).Value = 1.034430440380157
$ws.Cells.Item(7, 5).Value = 1.030091151219882
$ws.Cells.Item(7, 6).Value = 1.041682831265386
$ws.Cells.Item(7, 9).Value = 1.03190671428796
$ws.Cells.Item(7, 10).Value = 1.031126393106523
$ws.Cells.Item(7, 11).Value = 1.036917688588823
$ws.Cells.Item(7, 12).Value = 1.032589455744248
$ws.Cells.Item(7, 13).Value = 1.044151822199229
$ws.Cells.Item(7, 14).Value = 1.014409244169333

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.025320996108215
$ws.Cells.Item(8, 4).Value = 1.033489995048158
$ws.Cells.Item(8, 5).Value = 1.028929780573294
$ws.Cells.Item(8, 6).Value = 1.040489384278258
$ws.Cells.Item(8, 9).Value = 1.031694997612711
$ws.Cells.Item(8, 10).Value = 1.030370636169868
$ws.Cells.Item(8, 11).Value = 1.036228545087209
$ws.Cells.Item(8, 12).Value = 1.031681202525426
$ws.Cells.Item(8, 13).Value = 1.04320841541111
$ws.Cells.Item(8, 14).Value = 1.014157945447883

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.023155564307346
$ws.Cells.Item(9, 4).Value = 1.031831567101224
$ws.Cells.Item(9, 5).Value = 1.026888207722974
$ws.Cells.Item(9, 6).Value = 1.038389228909845
$ws.Cells.Item(9, 9).Value = 1.031307781151789
$ws.Cells.Item(9, 10).Value = 1.029036482261508
$ws.Cells.Item(9, 11).Value = 1.035007803073463
$ws.Cells.Item(9, 12).Value = 1.030080798678659
$ws.Cells.Item(9, 13).Value = 1.04154403236359
$ws.Cells.Item(9, 14).Value = 1.013714025377486

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.021713454493821
$ws.Cells.Item(10, 4).Value = 1.030725196851459
$ws.Cells.Item(10, 5).Value = 1.025530535821388
$ws.Cells.Item(10, 6).Value = 1.03699112425682
$ws.Cells.Item(10, 9).Value = 1.031040169763007
$ws.Cells.Item(10, 10).Value = 1.028145478658806
$ws.Cells.Item(10, 11).Value = 1.034189748147118
$ws.Cells.Item(10, 12).Value = 1.029013962388458
$ws.Cells.Item(10, 13).Value = 1.040433185394591
$ws.Cells.Item(10, 14).Value = 1.013417360011882

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.021089374032724
$ws.Cells.Item(11, 4).Value = 1.030245965438406
$ws.Cells.Item(11, 5).Value = 1.024943465234315
$ws.Cells.Item(11, 6).Value = 1.036386220533139
$ws.Cells.Item(11, 9).Value = 1.030922054414817
$ws.Cells.Item(11, 10).Value = 1.02775930371126
$ws.Cells.Item(11, 11).Value = 1.033834532055986
$ws.Cells.Item(11, 12).Value = 1.028552047053407
$ws.Cells.Item(11, 13).Value = 1.039951892894955
$ws.Cells.Item(11, 14).Value = 1.013288734427174

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.02085761801995
$ws.Cells.Item(12, 4).Value = 1.030067934006213
$ws.Cells.Item(12, 5).Value = 1.024725524022222
$ws.Cells.Item(12, 6).Value = 1.036161606510301
$ws.Cells.Item(12, 9).Value = 1.030877845566123
$ws.Cells.Item(12, 10).Value = 1.027615807243422
$ws.Cells.Item(12, 11).Value = 1.033702440833173
$ws.Cells.Item(12, 12).Value = 1.028380476910735
$ws.Cells.Item(12, 13).Value = 1.039773077252991
$ws.Cells.Item(12, 14).Value = 1.013240932317621

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020907327958346
$ws.Cells.Item(13, 4).Value = 1.030106123396556
$ws.Cells.Item(13, 5).Value = 1.024772267569295
$ws.Cells.Item(13, 6).Value = 1.036209783604001
$ws.Cells.Item(13, 9).Value = 1.030887343687958
$ws.Cells.Item(13, 10).Value = 1.027646590148014
$ws.Cells.Item(13, 11).Value = 1.033730781536617
$ws.Cells.Item(13, 12).Value = 1.0284172790112
$ws.Cells.Item(13, 13).Value = 1.039811435705338
$ws.Cells.Item(13, 14).Value = 1.013251187151209

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.021070215876657
$ws.Cells.Item(14, 4).Value = 1.030231249777201
$ws.Cells.Item(14, 5).Value = 1.024925447630778
$ws.Cells.Item(14, 6).Value = 1.036367652347645
$ws.Cells.Item(14, 9).Value = 1.03091840693888
$ws.Cells.Item(14, 10).Value = 1.027747443342971
$ws.Cells.Item(14, 11).Value = 1.033823616370632
$ws.Cells.Item(14, 12).Value = 1.028537864888224
$ws.Cells.Item(14, 13).Value = 1.039937112791963
$ws.Cells.Item(14, 14).Value = 1.013284783594728

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.021170583844637
$ws.Cells.Item(15, 4).Value = 1.030308341190319
$ws.Cells.Item(15, 5).Value = 1.025019843249378
$ws.Cells.Item(15, 6).Value = 1.036464930356249
$ws.Cells.Item(15, 9).Value = 1.030937501589198
$ws.Cells.Item(15, 10).Value = 1.027809575208801
$ws.Cells.Item(15, 11).Value = 1.033880795374925
$ws.Cells.Item(15, 12).Value = 1.028612162622577
$ws.Cells.Item(15, 13).Value = 1.040014541038385
$ws.Cells.Item(15, 14).Value = 1.013305480190614

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.021754880364933
$ws.Cells.Item(16, 4).Value = 1.030756998485273
$ws.Cells.Item(16, 5).Value = 1.025569514937829
$ws.Cells.Item(16, 6).Value = 1.037031280051661
$ws.Cells.Item(16, 9).Value = 1.031047961597364
$ws.Cells.Item(16, 10).Value = 1.028171100212403
$ws.Cells.Item(16, 11).Value = 1.034213301837173
$ws.Cells.Item(16, 12).Value = 1.029044618943698
$ws.Cells.Item(16, 13).Value = 1.040465121214508
$ws.Cells.Item(16, 14).Value = 1.013425892966831

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.022121491765703
$ws.Cells.Item(17, 4).Value = 1.031038385928284
$ws.Cells.Item(17, 5).Value = 1.025914527258664
$ws.Cells.Item(17, 6).Value = 1.037386666982101
$ws.Cells.Item(17, 9).Value = 1.03111665157503
$ws.Cells.Item(17, 10).Value = 1.028397778318556
$ws.Cells.Item(17, 11).Value = 1.034421609249931
$ws.Cells.Item(17, 12).Value = 1.02931589659371
$ws.Cells.Item(17, 13).Value = 1.040747681850776
$ws.Cells.Item(17, 14).Value = 1.013501380101326

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.022335364926956
$ws.Cells.Item(18, 4).Value = 1.031202498375485
$ws.Cells.Item(18, 5).Value = 1.026115845210438
$ws.Cells.Item(18, 6).Value = 1.037594004855846
$ws.Cells.Item(18, 9).Value = 1.031156501312514
$ws.Cells.Item(18, 10).Value = 1.028529960595207
$ws.Cells.Item(18, 11).Value = 1.034543015569821
$ws.Cells.Item(18, 12).Value = 1.029474131260516
$ws.Cells.Item(18, 13).Value = 1.040912466782671
$ws.Cells.Item(18, 14).Value = 1.013545394284187

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.022408296049594
$ws.Cells.Item(19, 4).Value = 1.031258453708688
$ws.Cells.Item(19, 5).Value = 1.026184502640278
$ws.Cells.Item(19, 6).Value = 1.037664709592036
$ws.Cells.Item(19, 9).Value = 1.031170052391667
$ws.Cells.Item(19, 10).Value = 1.028575025334069
$ws.Cells.Item(19, 11).Value = 1.034584395695124
$ws.Cells.Item(19, 12).Value = 1.0295280856717
$ws.Cells.Item(19, 13).Value = 1.040968649417338
$ws.Cells.Item(19, 14).Value = 1.01356039922776

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.022082154212909
$ws.Cells.Item(20, 4).Value = 1.03100819737281
$ws.Cells.Item(20, 5).Value = 1.025877502614583
$ws.Cells.Item(20, 6).Value = 1.037348532488902
$ws.Cells.Item(20, 9).Value = 1.031109304123888
$ws.Cells.Item(20, 10).Value = 1.028373461540485
$ws.Cells.Item(20, 11).Value = 1.034399269744758
$ws.Cells.Item(20, 12).Value = 1.02928679075131
$ws.Cells.Item(20, 13).Value = 1.040717368646752
$ws.Cells.Item(20, 14).Value = 1.013493282717961

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.021022247912936
$ws.Cells.Item(21, 4).Value = 1.030194403815925
$ws.Cells.Item(21, 5).Value = 1.024880336517093
$ws.Cells.Item(21, 6).Value = 1.036321161879876
$ws.Cells.Item(21, 9).Value = 1.030909268843217
$ws.Cells.Item(21, 10).Value = 1.027717746071666
$ws.Cells.Item(21, 11).Value = 1.033796282907353
$ws.Cells.Item(21, 12).Value = 1.02850235518541
$ws.Cells.Item(21, 13).Value = 1.039900105177433
$ws.Cells.Item(21, 14).Value = 1.013274890962566

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.020356162845092
$ws.Cells.Item(22, 4).Value = 1.029682603937157
$ws.Cells.Item(22, 5).Value = 1.024254090519935
$ws.Cells.Item(22, 6).Value = 1.035675642878807
$ws.Cells.Item(22, 9).Value = 1.030781557663122
$ws.Cells.Item(22, 10).Value = 1.027305160320939
$ws.Cells.Item(22, 11).Value = 1.033416304499376
$ws.Cells.Item(22, 12).Value = 1.028009183342188
$ws.Cells.Item(22, 13).Value = 1.039386015736765
$ws.Cells.Item(22, 14).Value = 1.013137435841101

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.020709236540733
$ws.Cells.Item(23, 4).Value = 1.029953931134222
$ws.Cells.Item(23, 5).Value = 1.024586007490304
$ws.Cells.Item(23, 6).Value = 1.036017803548494
$ws.Cells.Item(23, 9).Value = 1.030849443545996
$ws.Cells.Item(23, 10).Value = 1.027523909063703
$ws.Cells.Item(23, 11).Value = 1.033617819103277
$ws.Cells.Item(23, 12).Value = 1.028270619518258
$ws.Cells.Item(23, 13).Value = 1.039658566968748
$ws.Cells.Item(23, 14).Value = 1.013210316899982

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.022099929051825
$ws.Cells.Item(24, 4).Value = 1.031021838331404
$ws.Cells.Item(24, 5).Value = 1.025894232216307
$ws.Cells.Item(24, 6).Value = 1.037365763680734
$ws.Cells.Item(24, 9).Value = 1.031112624788129
$ws.Cells.Item(24, 10).Value = 1.028384449354694
$ws.Cells.Item(24, 11).Value = 1.034409364301847
$ws.Cells.Item(24, 12).Value = 1.029299942419435
$ws.Cells.Item(24, 13).Value = 1.040731065965144
$ws.Cells.Item(24, 14).Value = 1.013496941626991

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.023715117792803
$ws.Cells.Item(25, 4).Value = 1.032260449092179
$ws.Cells.Item(25, 5).Value = 1.027415413562096
$ws.Cells.Item(25, 6).Value = 1.038931823093133
$ws.Cells.Item(25, 9).Value = 1.031409557654283
$ws.Cells.Item(25, 10).Value = 1.029381673360419
$ws.Cells.Item(25, 11).Value = 1.035324143885474
$ws.Cells.Item(25, 12).Value = 1.030494528481268
$ws.Cells.Item(25, 13).Value = 1.041974542365786
$ws.Cells.Item(25, 14).Value = 1.013828917531855

Write-Host "Applied vm_pu.xlsx updates for Case_4_125 (380 kV case)"